$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Final" exam score for student 1 (H8) and student 3's "Multiplier" row (D10)
$ws.Range("H8").Value = 0.785
$ws.Range("D10").Value = 0.83

# Move the active selection to J13, as it was left when the grades were saved
$ws.Activate()
$ws.Range("J13").Select()
